$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 61-63: column A was empty, now carries the placeholder doc-number "N0000"
$ws.Range("A61").Value = "N0000"
$ws.Range("A62").Value = "N0000"
$ws.Range("A63").Value = "N0000"

# Rows 70-82: the CR_EWH checklist rows (columns A, B, E, F) are wiped out,
# leaving the row's styling in place but with empty cells.
$ws.Range("A70:B82").Value = ""
$ws.Range("E70:F82").Value = ""

# Move the worksheet selection to reflect the newly cleared block.
$ws.Range("A70:F82").Select()
